# Auto-generated edit script: applies numeric corrections to the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 50601
$ws.Range("J93").Value = 50601
$ws.Range("L93").Value = 50601
$ws.Range("N93").Value = -55593
$ws.Range("H95").Value = 31500
$ws.Range("J95").Value = 31500
$ws.Range("L95").Value = 31500
$ws.Range("N95").Value = -36992
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H116").Value = 5749.5
$ws.Range("I116").Value = 4999.5
$ws.Range("J116").Value = 6499.5
$ws.Range("K116").Value = 4999.5
$ws.Range("L116").Value = 6499.5
$ws.Range("M116").Value = -1557.5
$ws.Range("N116").Value = -13383.5
$ws.Range("H132").Value = 4267.6772
$ws.Range("I132").Value = 4387.069
$ws.Range("J132").Value = 2536.5
$ws.Range("K132").Value = 13161.207
$ws.Range("L132").Value = 7609.5
$ws.Range("M132").Value = -10631.207
$ws.Range("N132").Value = -12669.5
$ws.Range("H135").Value = 1141.8334
$ws.Range("I135").Value = 1141.8334
$ws.Range("K135").Value = 10276.5006
$ws.Range("M135").Value = -7741.500599999999
$ws.Range("H137").Value = 1898.25
$ws.Range("I137").Value = 2026.7858
$ws.Range("J137").Value = 998.5
$ws.Range("K137").Value = 6080.357400000001
$ws.Range("L137").Value = 2995.5
$ws.Range("M137").Value = -3530.357400000001
$ws.Range("N137").Value = -8095.5
$ws.Range("H138").Value = 1858.3077
$ws.Range("I138").Value = 640.7778
$ws.Range("K138").Value = 1922.3334
$ws.Range("M138").Value = 3217.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 533
$ws.Range("I5").Value = 616.8333
$ws.Range("K5").Value = 616.8333
$ws.Range("M5").Value = -504.8333
$ws.Range("H32").Value = 4102.387
$ws.Range("I32").Value = 2974.1072
$ws.Range("K32").Value = 2974.1072
$ws.Range("M32").Value = -2687.1072
$ws.Range("H37").Value = 5017
$ws.Range("I37").Value = 5017
$ws.Range("K37").Value = 5017
$ws.Range("M37").Value = -4744
$ws.Range("H80").Value = 63904.91
$ws.Range("J80").Value = 63904.91
$ws.Range("L80").Value = 63904.91
$ws.Range("N80").Value = -65900.91
$ws.Range("H83").Value = 63904.91
$ws.Range("J83").Value = 63904.91
$ws.Range("L83").Value = 191714.73
$ws.Range("N83").Value = -201698.73
$ws.Range("H94").Value = 45499.75
$ws.Range("J94").Value = 45499.75
$ws.Range("L94").Value = 45499.75
$ws.Range("N94").Value = -47301.75
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 533
$ws.Range("I4").Value = 616.8333
$ws.Range("K4").Value = 616.8333
$ws.Range("M4").Value = -501.8333
$ws.Range("H94").Value = 492.07693
$ws.Range("I94").Value = 434
$ws.Range("K94").Value = 434
$ws.Range("M94").Value = 17

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5714761.5
$ws.Range("I22").Value = 546.2
$ws.Range("K22").Value = 546.2
$ws.Range("M22").Value = -196.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 5449.5
$ws.Range("J32").Value = 9999
$ws.Range("L32").Value = 29997
$ws.Range("N32").Value = -30563
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H70").Value = 13100.9
$ws.Range("I70").Value = 5752.5
$ws.Range("J70").Value = 17999.834
$ws.Range("K70").Value = 17257.5
$ws.Range("L70").Value = 53999.50199999999
$ws.Range("M70").Value = -16942.5
$ws.Range("N70").Value = -54629.50199999999
$ws.Range("H73").Value = 13100.9
$ws.Range("I73").Value = 5752.5
$ws.Range("J73").Value = 17999.834
$ws.Range("K73").Value = 17257.5
$ws.Range("L73").Value = 53999.50199999999
$ws.Range("M73").Value = -16165.5
$ws.Range("N73").Value = -56183.50199999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.2
$ws.Range("I2").Value = 144
$ws.Range("K2").Value = 144
$ws.Range("M2").Value = -31
$ws.Range("H102").Value = 1835.6666
$ws.Range("I102").Value = 1753.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1753.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -131.5
$ws.Range("N102").Value = -5244
$ws.Range("H107").Value = 3487.6667
$ws.Range("I107").Value = 2319.4
$ws.Range("J107").Value = 4948
$ws.Range("K107").Value = 2319.4
$ws.Range("L107").Value = 4948
$ws.Range("M107").Value = -399.4000000000001
$ws.Range("N107").Value = -8788
$ws.Range("H122").Value = 2481.16
$ws.Range("I122").Value = 2006.1904
$ws.Range("J122").Value = 4974.75
$ws.Range("K122").Value = 6018.5712
$ws.Range("L122").Value = 14924.25
$ws.Range("M122").Value = -3568.5712
$ws.Range("N122").Value = -19824.25
$ws.Range("H126").Value = 1665.7
$ws.Range("I126").Value = 1665.7
$ws.Range("K126").Value = 4997.1
$ws.Range("M126").Value = -2527.1
$ws.Range("H132").Value = 1966.4
$ws.Range("I132").Value = 1775.5385
$ws.Range("K132").Value = 5326.6155
$ws.Range("M132").Value = -2796.6155

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1581.5
$ws.Range("I22").Value = 1442
$ws.Range("K22").Value = 1442
$ws.Range("M22").Value = -1147
$ws.Range("H27").Value = 1581.5
$ws.Range("I27").Value = 1442
$ws.Range("K27").Value = 1442
$ws.Range("M27").Value = -1335
$ws.Range("H46").Value = 1790
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 1222.3636
$ws.Range("I55").Value = 232.66667
$ws.Range("K55").Value = 232.66667
$ws.Range("M55").Value = -59.66667000000001
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 47108
$ws.Range("I45").Value = 36644.668
$ws.Range("J45").Value = 54955.5
$ws.Range("K45").Value = 36644.668
$ws.Range("L45").Value = 54955.5
$ws.Range("M45").Value = -36153.668
$ws.Range("N45").Value = -55937.5
$ws.Range("H101").Value = 22367.889
$ws.Range("J101").Value = 22367.889
$ws.Range("L101").Value = 22367.889
$ws.Range("N101").Value = -28857.889
$ws.Range("H107").Value = 13999
$ws.Range("I107").Value = 18000
$ws.Range("K107").Value = 54000
$ws.Range("M107").Value = -52080
$ws.Range("H122").Value = 2456.6287
$ws.Range("J122").Value = 2409.3333
$ws.Range("L122").Value = 7227.999899999999
$ws.Range("N122").Value = -12127.9999
$ws.Range("H132").Value = 3058.4119
$ws.Range("I132").Value = 3093.3125
$ws.Range("K132").Value = 9279.9375
$ws.Range("M132").Value = -6749.9375
$ws.Range("H136").Value = 7944.6665
$ws.Range("I136").Value = 5792.1665
$ws.Range("K136").Value = 17376.4995
$ws.Range("M136").Value = -14826.4995

